$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "Day 2 (28/05/2019)"
$ws.Range("B2").Value = "Complete First Steps with Sass Mixins Extends and Functions"

$ws.Columns.Item(2).ColumnWidth = 56.1666666666667

$ws.Range("B10").Select()
